# Updated symbol list on Sun Feb  5 13:17:49 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  2  = @{ D = "327.81";        E = "-0.95%" }
  3  = @{ D = "42.98";         E = "3.64%" }
  4  = @{ D = "5.578";         E = "-1.32%" }
  5  = @{ D = "0.08200";       E = "-1.70%" }
  6  = @{ D = "8.784";         E = "-0.28%" }
  7  = @{ D = "1.921";         E = "-5.09%" }
  8  = @{ D = "2.821";         E = "-4.52%" }
  9  = @{ D = "0.9476";        E = "2.11%" }
  10 = @{ D = "0.1218";        E = "-5.58%" }
  11 = @{ D = "0.1906";        E = "-2.72%" }
  12 = @{ D = "0.09811";       E = "4.68%" }
  13 = @{ D = "0.04581";       E = "18.66%" }
  14 = @{ D = "0.1070";        E = "0.42%" }
  15 = @{ D = "0.001292";      E = "-0.42%" }
  16 = @{ D = "0.006003";      E = "-1.83%" }
  17 = @{ D = "3.481";         E = "1.08%" }
  18 = @{ D = "4.514";         E = "-0.51%" }
  19 = @{ E = "0.03%" }
  20 = @{ D = "8.772";         E = "5.16%" }
  21 = @{ D = "0.1367";        E = "-0.30%" }
  22 = @{ D = "0.2731";        E = "11.31%" }
  23 = @{ D = "0.04427";       E = "0.80%" }
  24 = @{ D = "0.001249";      E = "-0.40%" }
  25 = @{ D = "0.004350";      E = "0.27%" }
  26 = @{ D = "0.0001240";     E = "3.56%" }
  27 = @{ D = "0.0004024";     E = "32.15%" }
  28 = @{ }
  29 = @{ }
  30 = @{ }
  31 = @{ }
  32 = @{ }
  33 = @{ }
  34 = @{ }
  35 = @{ }
  36 = @{ }
  37 = @{ }
  38 = @{ }
  39 = @{ D = "0.02734";       E = "-1.62%" }
  40 = @{ D = "0.05684";       E = "3.03%" }
  41 = @{ D = "0.007906";      E = "1.49%" }
  42 = @{ D = "0.009642";      E = "8.01%" }
  43 = @{ D = "0.1417";        E = "-1.37%" }
  44 = @{ D = "0.002159";      E = "-3.55%" }
  45 = @{ D = "0.009742";      E = "-13.51%" }
  46 = @{ D = "0.00007296";    E = "3.81%" }
  47 = @{ D = "0.00000000757"; E = "1.15%" }
  48 = @{ D = "0.003376";      E = "6.39%" }
  49 = @{ D = "0.002290";      E = "0.43%" }
  50 = @{ D = "0.00002119";    E = "1.15%" }
  51 = @{ D = "0.0002018";     E = "1.15%" }
}

foreach ($r in 2..51) {
  $info = $rows[$r]
  if ($info.ContainsKey("D")) {
    $ws.Range("D$r").NumberFormat = "@"
    $ws.Range("D$r").Value = $info["D"]
  }
  if ($info.ContainsKey("E")) {
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $info["E"]
  }
  $ws.Range("G$r").NumberFormat = "@"
  $ws.Range("G$r").Value = "13"
}
